$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Append new row to Logs sheet (row 9)
$logs.Range("A9").Value = "Wat zijn jullie openingstijden?"
$logs.Range("B9").Value = "mailmind.test@zohomail.eu"
$logs.Range("C9").Value = "Hallo, ik zou graag willen weten wat jullie openingstijden zijn. Dank je wel!"
$logs.Range("D9").Value = "Informatieaanvraag"
$logs.Range("E9").Value = "Beste klant,`nBedankt voor je bericht. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 18:00 uur. Op zaterdag zijn we geopend van 10:00 tot 15:00 uur. Op zondag zijn we gesloten. Mocht je verder nog vragen hebben, dan hoor ik graag van je.`nMet vriendelijke groet,`nE-mailassistent"
$logs.Range("F9").Value = "2025-06-19 11:58:12"
$logs.Range("G9").Value = "Ja"

# Append new row to Dashboard sheet (row 5)
$dash.Range("A5").Value = "Informatieaanvraag"
$dash.Range("B5").Value = 1

# Extend conditional formatting ranges on Logs sheet to include row 9
$dFcs = $logs.Range("D2:D8").FormatConditions
for ($i = 1; $i -le $dFcs.Count; $i++) {
    $dFcs.Item($i).ModifyAppliesToRange($logs.Range("D2:D9"))
}
$gFcs = $logs.Range("G2:G8").FormatConditions
for ($i = 1; $i -le $gFcs.Count; $i++) {
    $gFcs.Item($i).ModifyAppliesToRange($logs.Range("G2:G9"))
}

# Extend chart series ranges on the Dashboard chart to include row 5
$chart = $dash.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$5,'Dashboard'!`$B`$2:`$B`$5,1)"
